$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Strategy & Uncertainty")

# --- Cell content updates: append ".xml" to all original XML config filenames ---
$ws.Range("D4").Value2 = "Colombia_RPS_High.xml" + [char]10 + "Colombia_Nuclear_Zero.xml"
$ws.Range("E4").Value2 = "Strategy_1_High_RPS.xml" + [char]10 + "Strategy_1_High_Nuclear.xml"
$ws.Range("F4").Value2 = "Colombia_RPS_Low.xml" + [char]10 + "Colombia_Nuclear_Normal.xml"
$ws.Range("G4").Value2 = "Strategy_1_Low_RPS.xml" + [char]10 + "Strategy_1_Low_Nuclear.xml"
$ws.Range("D5").Value2 = "Colombia_Bldg_ShellApplianceEff_High.xml" + [char]10 + "Colombia_IndustrialEff_High.xml"
$ws.Range("E5").Value2 = "Strategy_2_High_BldEE.xml" + [char]10 + "Strategy_2_High_IndEE.xml"
$ws.Range("F5").Value2 = "Colombia_Bldg_ShellApplianceEff_Low.xml" + [char]10 + "Colombia_IndustrialEff_Low.xml"
$ws.Range("G5").Value2 = "Strategy_2_Low_BldEE.xml" + [char]10 + "Strategy_2_Low_IndEE.xml"
$ws.Range("D6").Value2 = "transportation_UCD_CORE_RapidEVsw_Colombia.xml"
$ws.Range("E6").Value2 = "Strategy_3_High_ElecTrans.xml"
$ws.Range("F6").Value2 = "transportation_UCD_CORE_ModEVsw.xml"
$ws.Range("G6").Value2 = "Strategy_3_Low_ElecTrans.xml"
$ws.Range("D7").Value2 = "Colombia_Public_Transport_High.xml"
$ws.Range("E7").Value2 = "Strategy_4_High_PublicTrans.xml"
$ws.Range("F7").Value2 = "Colombia_Public_Transport_Normal.xml"
$ws.Range("G7").Value2 = "Strategy_4_Low_PublicTrans.xml"
$ws.Range("D8").Value2 = "land_constraint_Colombia_10_afforestation.xml"
$ws.Range("E8").Value2 = "Strategy_5_High_AFOLU.xml"
$ws.Range("D9").Value2 = "Colombia_Low_Meat.xml"
$ws.Range("E9").Value2 = "Strategy_6_High_Meat.xml"
$ws.Range("D13").Value2 = "Colombia_GDP_High.xml" + [char]10 + "Colombia_Population_High.xml"
$ws.Range("E13").Value2 = "Uncertainty_1_High_GDP.xml" + [char]10 + "Uncertainty_1_High_Population.xml"
$ws.Range("F13").Value2 = "Colombia_GDP_Low.xml" + [char]10 + "Colombia_Population_Low.xml"
$ws.Range("G13").Value2 = "Uncertainty_1_Low_GDP.xml" + [char]10 + "Uncertainty_1_Low_Population.xml"
$ws.Range("D14").Value2 = "transportation_UCD_CORE_RapidEVcost_Colombia_noPubTrninterp.xml"
$ws.Range("E14").Value2 = "Uncertainty_2_High_EVCost.xml"
$ws.Range("F14").Value2 = "transportation_UCD_CORE_ModEVcost_Colombia_noPubTrninterp.xml"
$ws.Range("G14").Value2 = "Uncertainty_2_Low_EVCost.xml"
$ws.Range("E15").Value2 = "Uncertainty_3_High_RECostSolar.xml Uncertainty_3_High_RECostWind.xml"
$ws.Range("D16").Value2 = "Global_CCS_Cost_Normal.xml"
$ws.Range("E16").Value2 = "Uncertainty_4_High_CCSCost.xml"
$ws.Range("F16").Value2 = "Global_CCS_Cost_High.xml"
$ws.Range("G16").Value2 = "Uncertainty_4_Low_CCSCost.xml"
$ws.Range("D17").Value2 = "ag_prodchange_rcp2p6_gfdl_pdssat.xml" + [char]10 + "hydro_impacts_GFDL-ESM2M_rcp2p6.xml" + [char]10 + "runoff_impacts_GFDL-ESM2M_rcp2p6.xml"
$ws.Range("E17").Value2 = "Uncertainty_5_High_Ag.xml" + [char]10 + "Uncertainty_5_High_Hydro.xml" + [char]10 + "Uncertainty_5_High_Runoff.xml"
$ws.Range("F17").Value2 = "ag_prodchange_rcp2p6_hadgem2_pdssat.xml" + [char]10 + "hydro_impacts_HadGEM2-ES_rcp2p6.xml" + [char]10 + "runoff_impacts_HadGEM2-ES_rcp2p6.xml"
$ws.Range("G17").Value2 = "Uncertainty_5_Low_Ag.xml" + [char]10 + "Uncertainty_5_Low_Hydro.xml" + [char]10 + "Uncertainty_5_Low_Runoff.xml"
$ws.Range("D18").Value2 = "Global_ag_trade_HOV_CL_25.xml"
$ws.Range("E18").Value2 = "Uncertainty_6_High_HOV-CL.xml"

# --- Column width adjustments (content grew due to ".xml" suffixes) ---
$ws.Columns("D").ColumnWidth = 44.25
$ws.Columns("E").ColumnWidth = 31.5
$ws.Columns("F").ColumnWidth = 38.25
$ws.Columns("G").ColumnWidth = 29.25

# --- Row 14 shrinks since the EV-cost entries no longer have a leading blank line ---
$ws.Rows(14).RowHeight = 31.5

# --- Selection moved ---
$ws.Range("F23").Select() | Out-Null
